# Scheduled runner update: refresh computed leve profit figures (price/profit
# columns H:N) across the per-job Sheets, matching the latest market pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1129.7333
$ws.Range("I135").Value = 1082.24
$ws.Range("J135").Value = 1367.2
$ws.Range("K135").Value = 9740.16
$ws.Range("L135").Value = 12304.8
$ws.Range("M135").Value = -7205.16
$ws.Range("N135").Value = -17374.8
$ws.Range("H137").Value = 23810460
$ws.Range("I137").Value = 29412622
$ws.Range("J137").Value = 1270.375
$ws.Range("K137").Value = 88237866
$ws.Range("L137").Value = 3811.125
$ws.Range("M137").Value = -88235316
$ws.Range("N137").Value = -8911.125
$ws.Range("H138").Value = 4060789.2
$ws.Range("I138").Value = 1602249.1
$ws.Range("J138").Value = 4978155
$ws.Range("K138").Value = 4806747.300000001
$ws.Range("L138").Value = 14934465
$ws.Range("M138").Value = -4801607.300000001
$ws.Range("N138").Value = -14944745

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2444.0667
$ws.Range("I61").Value = 1819.1364
$ws.Range("J61").Value = 4162.625
$ws.Range("K61").Value = 1819.1364
$ws.Range("L61").Value = 4162.625
$ws.Range("M61").Value = -1607.1364
$ws.Range("N61").Value = -4586.625
$ws.Range("H74").Value = 8480.842000000001
$ws.Range("I74").Value = 3048
$ws.Range("J74").Value = 11650
$ws.Range("K74").Value = 3048
$ws.Range("L74").Value = 11650
$ws.Range("M74").Value = -2174
$ws.Range("N74").Value = -13398
$ws.Range("H77").Value = 8480.842000000001
$ws.Range("I77").Value = 3048
$ws.Range("J77").Value = 11650
$ws.Range("K77").Value = 15240
$ws.Range("L77").Value = 58250
$ws.Range("M77").Value = -10872
$ws.Range("N77").Value = -66986
$ws.Range("H132").Value = 2955.8
$ws.Range("I132").Value = 2928
$ws.Range("J132").Value = 3005.2222
$ws.Range("K132").Value = 8784
$ws.Range("L132").Value = 9015.6666
$ws.Range("M132").Value = -6254
$ws.Range("N132").Value = -14075.6666
$ws.Range("H136").Value = 2444.0667
$ws.Range("I136").Value = 1819.1364
$ws.Range("J136").Value = 4162.625
$ws.Range("K136").Value = 5457.4092
$ws.Range("L136").Value = 12487.875
$ws.Range("M136").Value = -2907.4092
$ws.Range("N136").Value = -17587.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 984.25
$ws.Range("I36").Value = 984.25
$ws.Range("K36").Value = 984.25
$ws.Range("M36").Value = -450.25
$ws.Range("H99").Value = 1968.3334
$ws.Range("I99").Value = 2113.3333
$ws.Range("J99").Value = 1871.6666
$ws.Range("K99").Value = 2113.3333
$ws.Range("L99").Value = 1871.6666
$ws.Range("M99").Value = -615.3332999999998
$ws.Range("N99").Value = -4867.6666
$ws.Range("H105").Value = 3117.8667
$ws.Range("I105").Value = 3016.85
$ws.Range("J105").Value = 3319.9
$ws.Range("K105").Value = 3016.85
$ws.Range("L105").Value = 3319.9
$ws.Range("M105").Value = -1269.85
$ws.Range("N105").Value = -6813.9
$ws.Range("H134").Value = 4070.6897
$ws.Range("I134").Value = 3090.75
$ws.Range("J134").Value = 4762.4116
$ws.Range("K134").Value = 9272.25
$ws.Range("L134").Value = 14287.2348
$ws.Range("M134").Value = -6737.25
$ws.Range("N134").Value = -19357.2348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1494.08
$ws.Range("I31").Value = 1016.2
$ws.Range("J31").Value = 2210.9
$ws.Range("K31").Value = 1016.2
$ws.Range("L31").Value = 2210.9
$ws.Range("M31").Value = -721.2
$ws.Range("N31").Value = -2800.9
$ws.Range("H34").Value = 1494.08
$ws.Range("I34").Value = 1016.2
$ws.Range("J34").Value = 2210.9
$ws.Range("K34").Value = 1016.2
$ws.Range("L34").Value = 2210.9
$ws.Range("M34").Value = -814.2
$ws.Range("N34").Value = -2614.9
$ws.Range("H132").Value = 1939.6316
$ws.Range("I132").Value = 1547.9348
$ws.Range("J132").Value = 3577.6365
$ws.Range("K132").Value = 4643.8044
$ws.Range("L132").Value = 10732.9095
$ws.Range("M132").Value = -2113.8044
$ws.Range("N132").Value = -15792.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1854.8
$ws.Range("I5").Value = 1213.3334
$ws.Range("K5").Value = 3640.0002
$ws.Range("M5").Value = -3528.0002
$ws.Range("H131").Value = 1669.94
$ws.Range("I131").Value = 487.14285
$ws.Range("J131").Value = 1862.4884
$ws.Range("K131").Value = 1461.42855
$ws.Range("L131").Value = 5587.4652
$ws.Range("M131").Value = 3578.57145
$ws.Range("N131").Value = -15667.4652
$ws.Range("H135").Value = 1854.8
$ws.Range("I135").Value = 1213.3334
$ws.Range("K135").Value = 10920.0006
$ws.Range("M135").Value = -8385.000599999999
$ws.Range("H139").Value = 1788.7142
$ws.Range("I139").Value = 1482.2222
$ws.Range("K139").Value = 4446.6666
$ws.Range("M139").Value = 693.3334000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2367.7222
$ws.Range("I80").Value = 2228.5454
$ws.Range("J80").Value = 2586.4285
$ws.Range("K80").Value = 2228.5454
$ws.Range("L80").Value = 2586.4285
$ws.Range("M80").Value = -1230.5454
$ws.Range("N80").Value = -4582.4285
$ws.Range("H83").Value = 2367.7222
$ws.Range("I83").Value = 2228.5454
$ws.Range("J83").Value = 2586.4285
$ws.Range("K83").Value = 11142.727
$ws.Range("L83").Value = 12932.1425
$ws.Range("M83").Value = -6150.726999999999
$ws.Range("N83").Value = -22916.1425
$ws.Range("H122").Value = 1853387.6
$ws.Range("I122").Value = 3704206
$ws.Range("J122").Value = 2569.3333
$ws.Range("K122").Value = 11112618
$ws.Range("L122").Value = 7707.999899999999
$ws.Range("M122").Value = -11110168
$ws.Range("N122").Value = -12607.9999
$ws.Range("H132").Value = 2855.4138
$ws.Range("I132").Value = 2481.907
$ws.Range("K132").Value = 7445.721
$ws.Range("M132").Value = -4915.721

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5288.6665
$ws.Range("I132").Value = 4763.9443
$ws.Range("J132").Value = 6338.1113
$ws.Range("K132").Value = 14291.8329
$ws.Range("L132").Value = 19014.3339
$ws.Range("M132").Value = -11761.8329
$ws.Range("N132").Value = -24074.3339
$ws.Range("H136").Value = 4346.436
$ws.Range("I136").Value = 2542.1333
$ws.Range("J136").Value = 10360.777
$ws.Range("K136").Value = 7626.3999
$ws.Range("L136").Value = 31082.331
$ws.Range("M136").Value = -5076.3999
$ws.Range("N136").Value = -36182.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13890867
$ws.Range("I132").Value = 16668348
$ws.Range("K132").Value = 50005044
$ws.Range("M132").Value = -50002514
$ws.Range("H136").Value = 6556817.5
$ws.Range("I136").Value = 7775985
$ws.Range("K136").Value = 23327955
$ws.Range("M136").Value = -23325405
